# edit.ps1 - Apply "Updated Source and Data Exchange page" changes
#
# Summary of changes (see task diff):
#   1. The cached text of every "datetimeFigureOut" date field (on the
#      slide master, all 11 slide layouts, and the notes master) changes
#      from "5/19/19" to "4/14/2021".
#   2. On each of the 7 "resource diagram" slides (4, 6, 8, 10, 12, 14,
#      16) the label "supportingOrganization" is renamed to "reporter"
#      and the label "patient" is renamed to "subject".
#   3. On 3 of those slides (6, 14, 16) the "patient"/"subject" textbox
#      is also nudged/resized slightly to better fit the new label.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update cached date-field text everywhere it appears: "5/19/19" ->
#    "4/14/2021" (slide master, notes master, and every slide layout).
# ---------------------------------------------------------------------

function Update-DateField($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -ceq "5/19/19") {
                $tr.Text = "4/14/2021"
            }
        }
    }
}

# Slide master
Update-DateField $p.SlideMaster

# Notes master
Update-DateField $p.NotesMaster

# Every slide layout ("custom layout") hanging off the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateField $layouts.Item($li)
}

# ---------------------------------------------------------------------
# 2 & 3. Rename labels + nudge boxes on the 7 resource-diagram slides.
# ---------------------------------------------------------------------

function Rename-Label($slide, [string]$oldText, [string]$newText) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -ceq $oldText) {
                $tr.Text = $newText
                return $shp
            }
        }
    }
    return $null
}

$diagramSlideNumbers = @(4, 6, 8, 10, 12, 14, 16)

foreach ($n in $diagramSlideNumbers) {
    $slide = $p.Slides.Item($n)

    # "supportingOrganization" -> "reporter"
    Rename-Label $slide "supportingOrganization" "reporter" | Out-Null

    # "patient" -> "subject" (keep the returned shape so slides that
    # also need a position/size tweak can use it directly)
    $subjectShape = Rename-Label $slide "patient" "subject"

    switch ($n) {
        6 {
            # off 2971393,2588462 ext 640000,261610 -> off 2977321,2499506 ext 724560,261610
            $subjectShape.Left = 234.43472440944882
            $subjectShape.Top = 196.81149606299212
            $subjectShape.Width = 57.051968503937005
            $subjectShape.Height = 20.599212598425197
        }
        14 {
            # off 3048000,3443657 ext 1143000,261610 -> off 3134557,3417527 ext 918895,261610
            $subjectShape.Left = 246.8155118110236
            $subjectShape.Top = 269.0966141732284
            $subjectShape.Width = 72.35393700787401
            $subjectShape.Height = 20.599212598425197
        }
        16 {
            # off 3048000,3443657 ext 1143000,261610 -> off 3028021,3391690 ext 918895,261610
            $subjectShape.Left = 238.4268503937008
            $subjectShape.Top = 267.06220472440947
            $subjectShape.Width = 72.35393700787401
            $subjectShape.Height = 20.599212598425197
        }
    }
}
